# Rewrite the LR-pairs data table (Vcam1 -> Itgad) with refreshed TPM-based values.
# One extra target cluster ("MuSCs") is now included for every sending cluster,
# expanding the table from 15 to 20 data rows (plus the header row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @('Sending cluster', 'Ligand symbol', 'Receptor symbol', 'Target cluster', 'Ligand-expressing cells', 'Ligand detection rate', 'Ligand average expression value', 'Ligand total expression value', 'Ligand derived specificity of average expression value', 'Ligand derived specificity of total expression value', 'Receptor-expressing cells', 'Receptor detection rate', 'Receptor average expression value', 'Receptor total expression value', 'Receptor derived specificity of average expression value', 'Receptor derived specificity of total expression value', 'Edge average expression weight', 'Edge total expression weight', 'Edge average expression derived specificity', 'Edge total expression derived specificity')
$row2 = @('ECs', 'Vcam1', 'Itgad', 'FAPs', [double]'3', [double]'1', [double]'17.93632866666666', [double]'53.808986', [double]'0.1226979812530711', [double]'0.1347750935001359', [double]'1', [double]'0.3333333333333333', [double]'0.02455766666666667', [double]'0.073673', [double]'0.03892688336861549', [double]'0.03893808583788527', [double]'0.4404743806197778', [double]'3.964269425578', [double]'0.004776250005802871', [double]'0.005247884159517303')
$row3 = @('ECs', 'Vcam1', 'Itgad', 'Inflammatory-Mac', [double]'3', [double]'1', [double]'17.93632866666666', [double]'53.808986', [double]'0.1226979812530711', [double]'0.1347750935001359', [double]'3', [double]'1', [double]'0.2696353333333334', [double]'0.8089060000000001', [double]'0.427404741468018', [double]'0.4275277410011866', [double]'4.83626795881289', [double]'43.526411629316', [double]'0.0524416989561166', [double]'0.0576200912673368')
$row4 = @('ECs', 'Vcam1', 'Itgad', 'MuSCs', [double]'3', [double]'1', [double]'17.93632866666666', [double]'53.808986', [double]'0.1226979812530711', [double]'0.1347750935001359', [double]'1', [double]'0.5', [double]'0.0005445', [double]'0.001089', [double]'0.000863098611195871', [double]'0.0005755646638179122', [double]'0.009766330958999998', [double]'0.05859798575399999', [double]'0.0001059004572160627', [double]'7.757178138143337E-05')
$row5 = @('ECs', 'Vcam1', 'Itgad', 'Resolving-Mac', [double]'3', [double]'1', [double]'17.93632866666666', [double]'53.808986', [double]'0.1226979812530711', [double]'0.1347750935001359', [double]'3', [double]'1', [double]'0.336129', [double]'1.008387', [double]'0.5328052765521707', [double]'0.5329586084971102', [double]'6.028920218397999', [double]'54.260281965582', [double]'0.06537413183393563', [double]'0.07182954629190033')
$row6 = @('FAPs', 'Vcam1', 'Itgad', 'FAPs', [double]'3', [double]'1', [double]'57.44330666666667', [double]'172.32992', [double]'0.3929554311523962', [double]'0.4316338739568692', [double]'1', [double]'0.3333333333333333', [double]'0.02455766666666667', [double]'0.073673', [double]'0.03892688336861549', [double]'0.03893808583788527', [double]'1.410673577351111', [double]'12.69606219616', [double]'0.01529653023753334', [double]'0.01680699683467152')
$row7 = @('FAPs', 'Vcam1', 'Itgad', 'Inflammatory-Mac', [double]'3', [double]'1', [double]'57.44330666666667', [double]'172.32992', [double]'0.3929554311523962', [double]'0.4316338739568692', [double]'3', [double]'1', [double]'0.2696353333333334', [double]'0.8089060000000001', [double]'0.427404741468018', [double]'0.4275277410011866', [double]'15.48874514083556', [double]'139.39870626752', [double]'0.1679510144601435', [double]'0.1845354550723712')
$row8 = @('FAPs', 'Vcam1', 'Itgad', 'MuSCs', [double]'3', [double]'1', [double]'57.44330666666667', [double]'172.32992', [double]'0.3929554311523962', [double]'0.4316338739568692', [double]'1', [double]'0.5', [double]'0.0005445', [double]'0.001089', [double]'0.000863098611195871', [double]'0.0005755646638179122', [double]'0.03127788048', [double]'0.18766728288', [double]'0.0003391592868895079', [double]'0.0002484332055564085')
$row9 = @('FAPs', 'Vcam1', 'Itgad', 'Resolving-Mac', [double]'3', [double]'1', [double]'57.44330666666667', [double]'172.32992', [double]'0.3929554311523962', [double]'0.4316338739568692', [double]'3', [double]'1', [double]'0.336129', [double]'1.008387', [double]'0.5328052765521707', [double]'0.5329586084971102', [double]'19.30836122656', [double]'173.77525103904', [double]'0.2093687271678299', [double]'0.23004298884427')
$row10 = @('Inflammatory-Mac', 'Vcam1', 'Itgad', 'FAPs', [double]'3', [double]'1', [double]'9.626273333333334', [double]'28.87882', [double]'0.06585095126993876', [double]'0.07233263354328205', [double]'1', [double]'0.3333333333333333', [double]'0.02455766666666667', [double]'0.073673', [double]'0.03892688336861549', [double]'0.03893808583788527', [double]'0.2363988117622222', [double]'2.12758930586', [double]'0.002563372299797289', [double]'0.002816494293788616')
$row11 = @('Inflammatory-Mac', 'Vcam1', 'Itgad', 'Inflammatory-Mac', [double]'3', [double]'1', [double]'9.626273333333334', [double]'28.87882', [double]'0.06585095126993876', [double]'0.07233263354328205', [double]'3', [double]'1', [double]'0.2696353333333334', [double]'0.8089060000000001', [double]'0.427404741468018', [double]'0.4275277410011866', [double]'2.595583418991112', [double]'23.36025077092', [double]'0.02814500880295123', [double]'0.03092420741942603')
$row12 = @('Inflammatory-Mac', 'Vcam1', 'Itgad', 'MuSCs', [double]'3', [double]'1', [double]'9.626273333333334', [double]'28.87882', [double]'0.06585095126993876', [double]'0.07233263354328205', [double]'1', [double]'0.5', [double]'0.0005445', [double]'0.001089', [double]'0.000863098611195871', [double]'0.0005755646638179122', [double]'0.00524150583', [double]'0.03144903498', [double]'5.683586458701112E-05', [double]'4.163210790840337E-05')
$row13 = @('Inflammatory-Mac', 'Vcam1', 'Itgad', 'Resolving-Mac', [double]'3', [double]'1', [double]'9.626273333333334', [double]'28.87882', [double]'0.06585095126993876', [double]'0.07233263354328205', [double]'3', [double]'1', [double]'0.336129', [double]'1.008387', [double]'0.5328052765521707', [double]'0.5329586084971102', [double]'3.23566962926', [double]'29.12102666334', [double]'0.03508573430260323', [double]'0.038550299722159')
$row14 = @('MuSCs', 'Vcam1', 'Itgad', 'FAPs', [double]'2', [double]'1', [double]'39.29803649999999', [double]'78.59607299999999', [double]'0.2688281328564436', [double]'0.1968591842135532', [double]'1', [double]'0.3333333333333333', [double]'0.02455766666666667', [double]'0.073673', [double]'0.03892688336861549', [double]'0.03893808583788527', [double]'0.9650680810215', [double]'5.790408486128999', [double]'0.01046464137390545', [double]'0.007665319812883402')
$row15 = @('MuSCs', 'Vcam1', 'Itgad', 'Inflammatory-Mac', [double]'2', [double]'1', [double]'39.29803649999999', [double]'78.59607299999999', [double]'0.2688281328564436', [double]'0.1968591842135532', [double]'3', [double]'1', [double]'0.2696353333333334', [double]'0.8089060000000001', [double]'0.427404741468018', [double]'0.4275277410011866', [double]'10.596139171023', [double]'63.576835026138', [double]'0.1148984186228383', [double]'0.08416276232215684')
$row16 = @('MuSCs', 'Vcam1', 'Itgad', 'MuSCs', [double]'2', [double]'1', [double]'39.29803649999999', [double]'78.59607299999999', [double]'0.2688281328564436', [double]'0.1968591842135532', [double]'1', [double]'0.5', [double]'0.0005445', [double]'0.001089', [double]'0.000863098611195871', [double]'0.0005755646638179122', [double]'0.02139778087425', [double]'0.08559112349699999', [double]'0.0002320251881187756', [double]'0.0001133051901813422')
$row17 = @('MuSCs', 'Vcam1', 'Itgad', 'Resolving-Mac', [double]'2', [double]'1', [double]'39.29803649999999', [double]'78.59607299999999', [double]'0.2688281328564436', [double]'0.1968591842135532', [double]'3', [double]'1', [double]'0.336129', [double]'1.008387', [double]'0.5328052765521707', [double]'0.5329586084971102', [double]'13.2092097107085', [double]'79.25525826425098', [double]'0.1432330476715811', [double]'0.1049177968883316')
$row18 = @('Resolving-Mac', 'Vcam1', 'Itgad', 'FAPs', [double]'3', [double]'1', [double]'21.87880766666666', [double]'65.63642299999999', [double]'0.1496675034681502', [double]'0.1643992147861598', [double]'1', [double]'0.3333333333333333', [double]'0.02455766666666667', [double]'0.073673', [double]'0.03892688336861549', [double]'0.03893808583788527', [double]'0.5372924657421111', [double]'4.835632191678999', [double]'0.005826089451576539', [double]'0.006401390737024428')
$row19 = @('Resolving-Mac', 'Vcam1', 'Itgad', 'Inflammatory-Mac', [double]'3', [double]'1', [double]'21.87880766666666', [double]'65.63642299999999', [double]'0.1496675034681502', [double]'0.1643992147861598', [double]'3', [double]'1', [double]'0.2696353333333334', [double]'0.8089060000000001', [double]'0.427404741468018', [double]'0.4275277410011866', [double]'5.899299598137556', [double]'53.093696383238', [double]'0.06396860062596844', [double]'0.07028522491989579')
$row20 = @('Resolving-Mac', 'Vcam1', 'Itgad', 'MuSCs', [double]'3', [double]'1', [double]'21.87880766666666', [double]'65.63642299999999', [double]'0.1496675034681502', [double]'0.1643992147861598', [double]'1', [double]'0.5', [double]'0.0005445', [double]'0.001089', [double]'0.000863098611195871', [double]'0.0005755646638179122', [double]'0.0119130107745', [double]'0.07147806464699999', [double]'0.0001291778143845137', [double]'9.462237879032483E-05')
$row21 = @('Resolving-Mac', 'Vcam1', 'Itgad', 'Resolving-Mac', [double]'3', [double]'1', [double]'21.87880766666666', [double]'65.63642299999999', [double]'0.1496675034681502', [double]'0.1643992147861598', [double]'3', [double]'1', [double]'0.336129', [double]'1.008387', [double]'0.5328052765521707', [double]'0.5329586084971102', [double]'7.354101742188998', [double]'66.18691567970099', [double]'0.07974363557622074', [double]'0.0876179767504493')

$rows = @($row1, $row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13, $row14, $row15, $row16, $row17, $row18, $row19, $row20, $row21)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value2 = $rowData[$j]
    }
}

